$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 82; existing rows 82-86 shift down to 83-87
$ws.Rows.Item(82).Insert()

# Populate the newly inserted row 82 with the new weekly price record
$ws.Range("A82").Value = 4
$ws.Range("B82").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C82").Value = "Los Lagos"
$ws.Range("D82").Value = 44747
$ws.Range("E82").Value = 10
$ws.Range("F82").Value = 100112026
$ws.Range("G82").Value = "Haba"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 80
$ws.Range("K82").Value = 20000
$ws.Range("L82").Value = 20000
$ws.Range("M82").Value = 20000
$ws.Range("N82").Value = "`$/saco 25 kilos"
$ws.Range("O82").Value = "Provincia de Limarí"
$ws.Range("P82").Value = 800
$ws.Range("Q82").Value = 25
$ws.Range("R82").Value = "Hortaliza"
